$d = $word.ActiveDocument

# 1. Merge "1." + "7" + " Assumptions" heading runs into a single run, without
#    touching the separate (bookmark-separated) " And Constraints" runs that
#    follow. Find/Replace on this paragraph tends to coalesce ALL runs with
#    matching formatting in the paragraph (even across the bookmark), so we
#    use a direct Range.Text assignment scoped tightly to "1.7 Assumptions"
#    instead - first overshooting by one char then trimming it, since a
#    no-op replacement (identical text) does not split/merge runs at all.
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "1.7 Assumptions And Constraints`r") {
        $headingPara = $d.Paragraphs($i)
        break
    }
}
$hStart = $headingPara.Range.Start
$rTarget = $d.Range($hStart, $hStart + 15)
$rTarget.Text = "1.7 AssumptionsX"
$rExtra = $d.Range($hStart + 15, $hStart + 16)
$rExtra.Text = ""

# 2. Merge the "The dataset acquired from Kaggle..." runs into a single run.
$d.Content.Find.Execute("The dataset acquired from Kaggle is a valid dataset, it is assumed that the records correspond to actual features and accurate measurements/description of the Indian dishes that were considered.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "The dataset acquired from Kaggle is a valid dataset, it is assumed that the records correspond to actual features and accurate measurements/description of the Indian dishes that were considered.", 2)

# 3. Merge the "All the features in the dataset..." runs into a single run.
$d.Content.Find.Execute("All the features in the dataset corresponding to " + [char]8220 + "indian_food" + [char]8221 + ", acquired from Kaggle are required for the analysis and the features are independent of each other and are collectively required for tackling the problem statement.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "All the features in the dataset corresponding to " + [char]8220 + "indian_food" + [char]8221 + ", acquired from Kaggle are required for the analysis and the features are independent of each other and are collectively required for tackling the problem statement.", 2)

# 4. Merge the "It is an assumption that the model is restricted..." runs into a single run.
$d.Content.Find.Execute("It is an assumption that the model is restricted to a binomial problem(Vegetarian and Non-Vegetarian) rather than multi class. Thus, this model can only be used to this specific problem statement.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "It is an assumption that the model is restricted to a binomial problem(Vegetarian and Non-Vegetarian) rather than multi class. Thus, this model can only be used to this specific problem statement.", 2)

# 5. Mark the built-in "Default Paragraph Font" character style as semi-hidden
#    (w:semiHidden), matching Word's normal "style list got re-saved" touch-up.
try {
    $defStyle = $d.Styles("Default Paragraph Font")
    $defStyle.Hidden = $true
} catch {
    # Some hosts do not expose a writable Style.Hidden; ignore if unsupported.
}
